# Append new scrape results as of 2025-10-01 01:53 JST.
#
# The scraper re-wrote the "取得日時" (fetched-at) timestamp on every existing
# row to the new run's time, and prepended a brand-new job posting
# ("CentOS server migration") above the previously-last row, pushing the
# old last row (VBA macro job) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-10-01 01:53:00"

# --- 1. Refresh the "取得日時" timestamp for every existing data row (2-16) ---
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# --- 2. Make room for the new row: insert a blank row above the old row 16 ---
#     (this shifts the old row 16 -- the VBA macro posting -- down to row 17,
#     carrying its hyperlink relationship with it conceptually; the engine
#     keeps the existing hyperlink anchored at F16 with the old target, so we
#     fix that target below and add a fresh hyperlink for the row that moved)
$ws.Rows.Item(16).Insert()

# --- 3. Fill the new row 16 with the new job posting ---
$ws.Range("A16").Value = $newTimestamp
$ws.Range("B16").Value = "【急募】CentOSサーバー移行の専門家を探しています!"
$ws.Range("C16").Value = "システム開発"
$ws.Range("D16").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E16").Value = "期限情報なし"
$ws.Range("F16").Value = "https://www.lancers.jp/work/detail/5404155"
$ws.Range("G16").Value = 13

# The pre-existing hyperlink (still anchored to F16 after the row insert)
# pointed at the old job's URL -- repoint it at the new job's URL.
$ws.Range("F16").Hyperlinks.Item(1).Address = "https://www.lancers.jp/work/detail/5404155"

# --- 4. The old "VBA macro" job, now on row 17, needs its own hyperlink ---
#     (row-insert did not carry a hyperlink down onto row 17 for us)
$ws.Hyperlinks.Add($ws.Range("F17"), "https://www.lancers.jp/work/detail/5404010")
